# Actualizacion desde MV -datos-
# Adds two new daily rows (04-10-2021 and 05-10-2021) to the interbank
# rates/amounts table, and fills in the remaining columns of the last
# existing row (01-10-2021) which previously only had the Serie + last
# column populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 190 (01-10-2021): fill in the previously-empty columns B:F ----
# (G190 already holds 1.5 in the source workbook.)
$ws.Range("B190").Value = 25000
$ws.Range("C190").Value = 1.5
$ws.Range("D190").Value = 1.5
$ws.Range("E190").Value = 1.5
$ws.Range("F190").Value = 3

# --- Row 191 (new): 04-10-2021 ------------------------------------------
# The "Serie" column stores the date as literal text (matching every
# other row in the sheet), not as a real Excel date. Typing the text
# directly via .Value would get auto-recognised as a date serial, so we
# build it as a text formula in a scratch cell, copy it, and paste just
# the value into place - that preserves the plain-text shared-string
# storage without pulling in a date number format.
$ws.Range("Z1").Formula = "=""04-10-2021"""
$ws.Range("Z1").Copy()
$ws.Range("A191").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("B191").Value = 95000
$ws.Range("C191").Value = 1.5
$ws.Range("D191").Value = 1.5
$ws.Range("E191").Value = 1.5
$ws.Range("F191").Value = 5
$ws.Range("G191").Value = 1.5

# --- Row 192 (new): 05-10-2021 ------------------------------------------
# Only Serie + last column are populated for this row.
$ws.Range("Z1").Formula = "=""05-10-2021"""
$ws.Range("Z1").Copy()
$ws.Range("A192").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("G192").Value = 1.5
